# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '29.091.52'
Set-TextValue $ws.Range('E2') '  -0.14%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.831.25'
Set-TextValue $ws.Range('E3') '  -0.75%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.9998'
Set-TextValue $ws.Range('E4') '  +0.20%  '

# Row 5
Set-TextValue $ws.Range('D5') '240.46'
Set-TextValue $ws.Range('E5') '  -2.13%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.6841'
Set-TextValue $ws.Range('E6') '  -1.96%  '

# Row 7
Set-TextValue $ws.Range('D7') '1.000'
Set-TextValue $ws.Range('E7') '  +0.12%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.3011'
Set-TextValue $ws.Range('E8') '  -1.58%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.07453'
Set-TextValue $ws.Range('E9') '  -3.38%  '

# Row 10
Set-TextValue $ws.Range('D10') '23.06'
Set-TextValue $ws.Range('E10') '  -2.24%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.07657'
Set-TextValue $ws.Range('E11') '  -2.17%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D12') '5.056'
Set-TextValue $ws.Range('E12') '  -1.37%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.825.30'
Set-TextValue $ws.Range('E13') '  -0.90%  '

# Row 14
Set-TextValue $ws.Range('D14') '0.6807'
Set-TextValue $ws.Range('E14') '  -0.68%  '

# Row 15
Set-TextValue $ws.Range('D15') '87.54'
Set-TextValue $ws.Range('E15') '  -5.88%  '

# Row 16
Set-TextValue $ws.Range('D16') '6.132'
Set-TextValue $ws.Range('E16') '  -7.84%  '

# Row 17
Set-TextValue $ws.Range('D17') '29.108.15'
Set-TextValue $ws.Range('E17') '  +0.08%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.000008184'
Set-TextValue $ws.Range('E18') '  -1.41%  '

# Row 19
Set-TextValue $ws.Range('D19') '2.084.81'
Set-TextValue $ws.Range('E19') '  +0.30%  '

# Row 20
Set-TextValue $ws.Range('D20') '227.91'
Set-TextValue $ws.Range('E20') '  -5.94%  '

# Row 21
Set-TextValue $ws.Range('D21') '12.50'
Set-TextValue $ws.Range('E21') '  -2.08%  '

# Row 22
Set-TextValue $ws.Range('D22') '1.000'
Set-TextValue $ws.Range('E22') '  +0.11%  '

# Row 23
Set-TextValue $ws.Range('D23') '7.418'
Set-TextValue $ws.Range('E23') '  -0.94%  '

# Row 24
Set-TextValue $ws.Range('D24') '1.000'
Set-TextValue $ws.Range('E24') '  +0.15%  '

# Row 25
Set-TextValue $ws.Range('D25') '0.1455'
Set-TextValue $ws.Range('E25') '  -3.59%  '

# Row 26
Set-TextValue $ws.Range('D26') '160.16'
Set-TextValue $ws.Range('E26') '  +0.58%  '

# Row 27
Set-TextValue $ws.Range('D27') '8.730'
Set-TextValue $ws.Range('E27') '  -0.96%  '

# Row 28
Set-TextValue $ws.Range('E28') '  -0.74%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.508'
Set-TextValue $ws.Range('E29') '  -2.15%  '

# Row 30
Set-TextValue $ws.Range('D30') '4.287'
Set-TextValue $ws.Range('E30') '  +1.25%  '

# Row 31
Set-TextValue $ws.Range('D31') '4.146'
Set-TextValue $ws.Range('E31') '  -0.78%  '

# Row 32
Set-TextValue $ws.Range('D32') '1.197'
Set-TextValue $ws.Range('E32') '  -0.37%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.05160'
Set-TextValue $ws.Range('E33') '  +0.86%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.7664'
Set-TextValue $ws.Range('E34') '  -2.51%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.840'
Set-TextValue $ws.Range('E35') '  -1.43%  '

# Row 36
Set-TextValue $ws.Range('D36') '1.131'
Set-TextValue $ws.Range('E36') '  -1.46%  '

# Row 37
Set-TextValue $ws.Range('D37') '2.674'
Set-TextValue $ws.Range('E37') '  -0.62%  '

# Row 38
Set-TextValue $ws.Range('D38') '1.304.63'
Set-TextValue $ws.Range('E38') '  -0.62%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.01839'
Set-TextValue $ws.Range('E39') '  -1.42%  '

# Row 40
Set-TextValue $ws.Range('D40') '2.718'
Set-TextValue $ws.Range('E40') '  +0.41%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.9298'
Set-TextValue $ws.Range('E41') '  -2.12%  '

# Row 42
Set-TextValue $ws.Range('D42') '5.824'
Set-TextValue $ws.Range('E42') '  -4.82%  '

# Row 43
Set-TextValue $ws.Range('D43') '104.56'
Set-TextValue $ws.Range('E43') '  -2.97%  '

# Row 44
Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  +0.10%  '

# Row 45
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D45') '1.986.85'
Set-TextValue $ws.Range('E45') '  +0.38%  '

# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.00000000123'
Set-TextValue $ws.Range('E46') '  +2.91%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D47') '65.05'
Set-TextValue $ws.Range('E47') '  +1.20%  '

# Row 48
Set-TextValue $ws.Range('D48') '0.5194'
Set-TextValue $ws.Range('E48') '  +0.46%  '

# Row 49
Set-TextValue $ws.Range('D49') '9.550'
Set-TextValue $ws.Range('E49') '  -1.76%  '

# Row 50
Set-TextValue $ws.Range('E50') '  +0.41%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.05923'
Set-TextValue $ws.Range('E51') '  +0.69%  '
